$wb = $excel.ActiveWorkbook

# Rename the worksheets to match the new business-context names.
$wsInternet = $wb.Worksheets.Item(1)
$wsInternet.Name = "Internet_Connection"

$wsOnline = $wb.Worksheets.Item(2)
$wsOnline.Name = "Online Purchase"

# Update the selected cell on each sheet.
$wsInternet.Range("C32").Select()

$wsOnline.Range("F19").Select()

# Make "Online Purchase" the active (visible/selected) sheet/tab.
$wsOnline.Activate()
$wsOnline.Range("F19").Select()
